$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ComparacionEscalado")
$ws.Range("I6").ClearContents()
$excel.Calculate()
Write-Output ("K6 after: " + $ws.Range("K6").Value)
